$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '27.007.70'
$ws.Cells.Item(2, 5).Value = '  +0.54%  '
$ws.Cells.Item(3, 4).Value = '1.561.74'
$ws.Cells.Item(3, 5).Value = '  +0.73%  '
$ws.Cells.Item(4, 5).Value = '  -0.15%  '
$ws.Cells.Item(5, 4).Value = "'207.51"
$ws.Cells.Item(5, 5).Value = '  +0.42%  '
$ws.Cells.Item(6, 5).Value = '  +1.19%  '
$ws.Cells.Item(8, 4).Value = "'22.15"
$ws.Cells.Item(8, 5).Value = '  +2.19%  '
$ws.Cells.Item(9, 5).Value = '  +0.15%  '
$ws.Cells.Item(10, 4).Value = "'0.0596"
$ws.Cells.Item(10, 5).Value = '  +1.89%  '
$ws.Cells.Item(11, 4).Value = "'0.0859"
$ws.Cells.Item(11, 5).Value = '  +0.07%  '
$ws.Cells.Item(12, 4).Value = '1.784.79'
$ws.Cells.Item(12, 5).Value = '  +0.68%  '
$ws.Cells.Item(13, 4).Value = '1.541.90'
$ws.Cells.Item(13, 5).Value = '  -0.82%  '
$ws.Cells.Item(15, 5).Value = '  +1.40%  '
$ws.Cells.Item(16, 4).Value = "'62.08"
$ws.Cells.Item(16, 5).Value = '  +0.68%  '
$ws.Cells.Item(17, 4).Value = '27.004.97'
$ws.Cells.Item(17, 5).Value = '  +0.54%  '
$ws.Cells.Item(18, 4).Value = '0.0₃0706'
$ws.Cells.Item(18, 5).Value = '  +2.61%  '
$ws.Cells.Item(19, 4).Value = "'217.15"
$ws.Cells.Item(19, 5).Value = '  +0.30%  '
$ws.Cells.Item(20, 5).Value = '  +2.39%  '
$ws.Cells.Item(21, 5).Value = '  -0.19%  '
$ws.Cells.Item(22, 5).Value = '  +1.57%  '
$ws.Cells.Item(23, 4).Value = "'9.25"
$ws.Cells.Item(23, 5).Value = '  +0.97%  '
$ws.Cells.Item(24, 5).Value = '  -2.67%  '
$ws.Cells.Item(25, 4).Value = "'153.19"
$ws.Cells.Item(25, 5).Value = '  -0.27%  '
$ws.Cells.Item(26, 5).Value = '  +0.23%  '
$ws.Cells.Item(27, 5).Value = '  +1.24%  '
$ws.Cells.Item(28, 4).Value = "'0.104"
$ws.Cells.Item(28, 5).Value = '  +1.50%  '
$ws.Cells.Item(29, 5).Value = '  -0.16%  '
$ws.Cells.Item(30, 5).Value = '  +1.18%  '
$ws.Cells.Item(31, 4).Value = "'1.12"
$ws.Cells.Item(31, 5).Value = '  +2.15%  '
$ws.Cells.Item(32, 4).Value = "'3.24"
$ws.Cells.Item(32, 5).Value = '  +0.63%  '
$ws.Cells.Item(33, 4).Value = '1.423.19'
$ws.Cells.Item(33, 5).Value = '  +0.48%  '
$ws.Cells.Item(34, 4).Value = "'3.11"
$ws.Cells.Item(34, 5).Value = '  +3.69%  '
$ws.Cells.Item(35, 5).Value = '  +3.09%  '
$ws.Cells.Item(36, 5).Value = '  +9.50%  '
$ws.Cells.Item(37, 5).Value = '  +1.20%  '
$ws.Cells.Item(38, 5).Value = '  +0.81%  '
$ws.Cells.Item(39, 4).Value = "'0.531"
$ws.Cells.Item(39, 5).Value = '  +1.85%  '
$ws.Cells.Item(40, 5).Value = '  +0.46%  '
$ws.Cells.Item(41, 5).Value = '  -0.17%  '
$ws.Cells.Item(42, 2).Value = 'FraxShare'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(42, 4).Value = "'5.67"
$ws.Cells.Item(42, 5).Value = '  +0.57%  '
$ws.Cells.Item(43, 2).Value = 'MXToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(43, 4).Value = "'2.33"
$ws.Cells.Item(43, 5).Value = '  +3.13%  '
$ws.Cells.Item(44, 4).Value = "'0.999"
$ws.Cells.Item(44, 5).Value = '  +0.99%  '
$ws.Cells.Item(45, 4).Value = "'64.96"
$ws.Cells.Item(45, 5).Value = '  +2.22%  '
$ws.Cells.Item(46, 5).Value = '  +0.41%  '
$ws.Cells.Item(47, 4).Value = '1.703.43'
$ws.Cells.Item(47, 5).Value = '  +0.97%  '
$ws.Cells.Item(48, 4).Value = "'87.61"
$ws.Cells.Item(48, 5).Value = '  +1.79%  '
$ws.Cells.Item(49, 5).Value = '  +0.80%  '
$ws.Cells.Item(50, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(50, 4).Value = '0.0₆0101'
$ws.Cells.Item(50, 5).Value = '  +0.39%  '
$ws.Cells.Item(51, 2).Value = 'Algorand'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(51, 4).Value = "'0.0957"
$ws.Cells.Item(51, 5).Value = '  -0.35%  '
